$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked_lbl" (D) and "is_enabled_lbl" (E) columns so the
# remaining columns (order_by, rem) shift left into D and E.
$ws.Range("D1:E1").EntireColumn.Delete()
